$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
# Columns E..I already existed in the source sheet, so they keep the bold/
# bordered/centered header style ("s=1") automatically when only .Value is
# reassigned. Columns J..M are brand-new cells with no style yet, so the
# header formatting is cloned from an existing header cell (format-only
# paste) before the new header text is written.
$ws.Range("E1").Value = "4. Agriculture land area (% of land area)_x"
$ws.Range("F1").Value = "4. Agriculture land area (% of land area)_y"
$ws.Range("G1").Value = "5. Average precipitation (mm)"
$ws.Range("H1").Value = "7. Fertilizer consumption (kilograms per hectare of arable land)_x"
$ws.Range("I1").Value = "7. Fertilizer consumption (kilograms per hectare of arable land)_y"

$ws.Range("A1").Copy()
$ws.Range("J1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("J1").Value = "13. Population_x"
$ws.Range("K1").Value = "13. Population_y"
$ws.Range("L1").Value = "17. Employment in agriculture (% of total employment) (modeled ILO estimate)_x"
$ws.Range("M1").Value = "17. Employment in agriculture (% of total employment) (modeled ILO estimate)_y"

# ---- Data rows (row 2 .. row 10) ----
# code, year, crop(text), temp, agri_x, agri_y, precip, fert_x, fert_y, pop_x, pop_y, emp_x, emp_y
$data = @(
    @("LBR", 2002, "84.59", 25.75, 14.82114099, 14.82114099, 2297.57, 7.79275, 7.79275, 3060599, 3060599, 51.2976220963373, 51.2976220963373),
    @("LBR", 2003, "84.04", 25.76, 15.13546408, 15.13546408, 2348.52, 8.677250000000001, 8.677250000000001, 3085173, 3085173, 53.021896148668, 53.021896148668),
    @("LBR", 2006, "87.51", 25.79, 16.06935631, 16.06935631, 2387.49, 11.685, 11.685, 3455397, 3455397, 50.7591451193298, 50.7591451193298),
    @("LBR", 2008, "97.42", 25.8, 16.68627284, 16.68627284, 2283.36, 7.82372093, 7.82372093, 3783887, 3783887, 49.1975944159694, 49.1975944159694),
    @("LBR", 2009, "90.28", 25.83, 16.99679298, 16.99679298, 2318.51, 7.827659574, 7.827659574, 3905066, 3905066, 48.2323209698238, 48.2323209698238),
    @("LBR", 2010, "91.79", 25.86, 17.30725291, 17.30725291, 2309.14, 12.41625, 12.41625, 4019956, 4019956, 47.3061268582341, 47.3061268582341),
    @("LBR", 2012, "99.98", 25.94, 17.92604028, 17.92604028, 2362.24, 11.4594, 11.4594, 4331740, 4331740, 44.8846559111217, 44.8846559111217),
    @("LBR", 2013, "95.93", 25.94, 18.23662064, 18.23662064, 2297.17, 11.9564, 11.9564, 4427313, 4427313, 43.4367755988557, 43.4367755988557),
    @("LBR", 2014, "90.47", 25.97, 18.54623443, 18.54623443, 2306.41, 12.0842, 12.0842, 4519398, 4519398, 42.8104478161898, 42.8104478161898)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    # Column C must stay a text string like "84.59" (not auto-converted to a number).
    # Leading apostrophe forces text entry; resetting the Style afterwards strips the
    # "quote prefix" formatting that Excel would otherwise stamp onto the cell.
    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $ws.Cells.Item($r, 3).Style = "Normal"

    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $r = $r + 1
}

Write-Output "done"
